$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text number format on Price/Volume cells so exact literal
# representations (trailing zeros, percent signs, etc.) are preserved
# exactly as authored, matching the original inlineStr cell content.
$numericCells = @("E2", "D3", "E3", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "D8", "E8", "D9", "E9", "D10", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "D18", "E18", "E19", "D20", "E20", "E21", "D22", "E22", "E23", "D24", "E24", "D25", "E25", "D26", "E26", "D27", "E27", "D39", "E39", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "D46", "E46", "D47", "E47", "D48", "E48", "E49", "D50", "E50", "D51", "E51")
foreach ($addr in $numericCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values cell by cell, matching the refreshed
# coinranking.com symbol list snapshot.
$ws.Range('E2').Value = '1.78%'
$ws.Range('D3').Value = '37.65'
$ws.Range('E3').Value = '0.94%'
$ws.Range('E4').Value = '0.55%'
$ws.Range('D5').Value = '0.07915'
$ws.Range('E5').Value = '1.96%'
$ws.Range('B6').Value = 'FTXToken'
$ws.Range('C6').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D6').Value = '1.920'
$ws.Range('E6').Value = '-0.25%'
$ws.Range('B7').Value = 'KuCoinToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range('D7').Value = '8.272'
$ws.Range('E7').Value = '0.53%'
$ws.Range('B8').Value = 'BTSEToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D8').Value = '2.921'
$ws.Range('E8').Value = '-7.33%'
$ws.Range('B9').Value = 'MXToken'
$ws.Range('C9').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D9').Value = '0.9216'
$ws.Range('E9').Value = '-0.03%'
$ws.Range('B10').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C10').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D10').Value = '0.1233'
$ws.Range('E10').Value = '-3.37%'
$ws.Range('B11').Value = 'WazirX'
$ws.Range('C11').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D11').Value = '0.1918'
$ws.Range('E11').Value = '0.98%'
$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D12').Value = '0.09153'
$ws.Range('E12').Value = '4.05%'
$ws.Range('B13').Value = 'BitrueCoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D13').Value = '0.03322'
$ws.Range('E13').Value = '-3.45%'
$ws.Range('B14').Value = 'BitMartToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D14').Value = '0.09601'
$ws.Range('E14').Value = '-1.30%'
$ws.Range('B15').Value = 'BitForexToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D15').Value = '0.001384'
$ws.Range('E15').Value = '1.23%'
$ws.Range('B16').Value = 'TigerCash'
$ws.Range('C16').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D16').Value = '0.005711'
$ws.Range('E16').Value = '-6.81%'
$ws.Range('B17').Value = 'LEO'
$ws.Range('C17').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D17').Value = '3.521'
$ws.Range('E17').Value = '-0.97%'
$ws.Range('B18').Value = 'GateToken'
$ws.Range('C18').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D18').Value = '4.422'
$ws.Range('E18').Value = '0.52%'
$ws.Range('E19').Value = '2.07%'
$ws.Range('D20').Value = '5.255'
$ws.Range('E20').Value = '4.37%'
$ws.Range('E21').Value = '-1.26%'
$ws.Range('D22').Value = '0.2592'
$ws.Range('E22').Value = '3.58%'
$ws.Range('E23').Value = '-0.71%'
$ws.Range('D24').Value = '0.04367'
$ws.Range('E24').Value = '0.48%'
$ws.Range('D25').Value = '0.001249'
$ws.Range('E25').Value = '1.99%'
$ws.Range('D26').Value = '0.004297'
$ws.Range('E26').Value = '-4.46%'
$ws.Range('D27').Value = '0.0001220'
$ws.Range('E27').Value = '-10.10%'
$ws.Range('D39').Value = '0.02224'
$ws.Range('E39').Value = '2.46%'
$ws.Range('D40').Value = '0.05127'
$ws.Range('E40').Value = '3.66%'
$ws.Range('D41').Value = '0.007438'
$ws.Range('E41').Value = '-3.28%'
$ws.Range('D42').Value = '0.1361'
$ws.Range('E42').Value = '1.74%'
$ws.Range('D43').Value = '0.008744'
$ws.Range('E43').Value = '-11.12%'
$ws.Range('D44').Value = '0.002010'
$ws.Range('E44').Value = '0.27%'
$ws.Range('D45').Value = '0.008626'
$ws.Range('E45').Value = '3.05%'
$ws.Range('D46').Value = '0.00006731'
$ws.Range('E46').Value = '-1.90%'
$ws.Range('D47').Value = '0.00000000750'
$ws.Range('E47').Value = '-0.52%'
$ws.Range('D48').Value = '0.003351'
$ws.Range('E48').Value = '11.07%'
$ws.Range('E49').Value = '-8.15%'
$ws.Range('D50').Value = '0.00002100'
$ws.Range('E50').Value = '-0.52%'
$ws.Range('D51').Value = '0.0002000'
$ws.Range('E51').Value = '-0.52%'
